# Actualización automática: rollover mensual de "LOZANO MOLINA TITO"
#
# - Hoja "VENTAS POR GRUPO": los importes del mes en curso (ya volcados a la
#   hoja "VENTA MENSUAL") se limpian a 0, junto con sus contadores "N de 30"
#   de la fila 32 que pasan a "0 de 30".
# - Hoja "VENTA MENSUAL": los meses se recorren una columna a la izquierda
#   (C<-D, D<-E, E<-F) y la nueva columna F queda en 0 a la espera del nuevo
#   mes; los encabezados de mes (fila 1) se recorren igual y el nuevo mes
#   ("octubre") se escribe en F1. Los anchos de columna de E y F se
#   intercambian para conservar el ancho "ancho" en la columna que ahora
#   aloja el nombre de mes más largo.

$wb = $excel.ActiveWorkbook

$wsGrupo   = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# --- Hoja "VENTAS POR GRUPO": limpiar a 0 los importes ya trasladados -----
$wsGrupo.Range("D2").Value  = 0
$wsGrupo.Range("L4").Value  = 0
$wsGrupo.Range("L6").Value  = 0
$wsGrupo.Range("O15").Value = 0
$wsGrupo.Range("M16").Value = 0
$wsGrupo.Range("C17").Value = 0
$wsGrupo.Range("L17").Value = 0
$wsGrupo.Range("L19").Value = 0
$wsGrupo.Range("L22").Value = 0
$wsGrupo.Range("D24").Value = 0
$wsGrupo.Range("H24").Value = 0
$wsGrupo.Range("I24").Value = 0

# Contadores "N de 30" (fila 32) de las columnas afectadas -> "0 de 30"
$wsGrupo.Range("C32").Value = "0 de 30"
$wsGrupo.Range("D32").Value = "0 de 30"
$wsGrupo.Range("H32").Value = "0 de 30"
$wsGrupo.Range("I32").Value = "0 de 30"
$wsGrupo.Range("L32").Value = "0 de 30"
$wsGrupo.Range("M32").Value = "0 de 30"
$wsGrupo.Range("O32").Value = "0 de 30"

# --- Hoja "VENTA MENSUAL": recorrer columnas de mes una posición a la ------
# --- izquierda (C<-D, D<-E, E<-F) dejando F en blanco (0) ------------------

# Encabezados de mes (fila 1)
$wsMensual.Range("C1").Value = "julio"
$wsMensual.Range("D1").Value = "agosto"
$wsMensual.Range("E1").Value = "septiembre"
$wsMensual.Range("F1").Value = "octubre"

# Anchos de columna: E y F intercambian su ancho.
# Nota: la propiedad COM `ColumnWidth` añade un relleno fijo de 5/6 de
# caracter frente al atributo crudo `width` del XML, así que se resta aquí
# para que el XML resultante quede en los valores exactos 16 y 13.
$wsMensual.Columns.Item(5).ColumnWidth = 16 - (5 / 6)
$wsMensual.Columns.Item(6).ColumnWidth = 13 - (5 / 6)

$lastRow = 32
for ($r = 2; $r -le $lastRow; $r++) {
    $d = $wsMensual.Cells.Item($r, 4).Value2
    $e = $wsMensual.Cells.Item($r, 5).Value2
    $f = $wsMensual.Cells.Item($r, 6).Value2

    $wsMensual.Cells.Item($r, 3).Value = $d
    $wsMensual.Cells.Item($r, 4).Value = $e
    $wsMensual.Cells.Item($r, 5).Value = $f
    $wsMensual.Cells.Item($r, 6).Value = 0
}
